$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each (cell, new value) pair below reproduces one <c> text change from the
# source diff. D/E columns store plain text that merely *looks* numeric
# (prices / percentages), so NumberFormat is forced to "@" (Text) right before
# the assignment -- otherwise Excel's normal type inference would silently
# coerce the literal string into a real number/percentage and mangle the
# exact text (trailing zeros, "%", precision, etc.).
$textCells = @(
    @{Cell = "D2"; Value = "246.79"}
    @{Cell = "E2"; Value = "0.43%"}
    @{Cell = "D3"; Value = "26.04"}
    @{Cell = "E3"; Value = "3.08%"}
    @{Cell = "D4"; Value = "5.199"}
    @{Cell = "E4"; Value = "2.95%"}
    @{Cell = "D5"; Value = "0.05594"}
    @{Cell = "E5"; Value = "-0.04%"}
    @{Cell = "D6"; Value = "6.483"}
    @{Cell = "E6"; Value = "-1.26%"}
    @{Cell = "D7"; Value = "0.8127"}
    @{Cell = "E7"; Value = "-0.50%"}
    @{Cell = "D8"; Value = "0.8458"}
    @{Cell = "E8"; Value = "1.36%"}
    @{Cell = "D9"; Value = "0.06918"}
    @{Cell = "E9"; Value = "-0.46%"}
    @{Cell = "D10"; Value = "0.02825"}
    @{Cell = "E10"; Value = "-0.07%"}
    @{Cell = "D11"; Value = "0.09381"}
    @{Cell = "E11"; Value = "-0.07%"}
    @{Cell = "D12"; Value = "0.001510"}
    @{Cell = "E12"; Value = "-0.67%"}
    @{Cell = "D13"; Value = "0.0005957"}
    @{Cell = "E13"; Value = "-0.15%"}
    @{Cell = "D14"; Value = "0.006173"}
    @{Cell = "E14"; Value = "0.06%"}
    @{Cell = "D15"; Value = "3.609"}
    @{Cell = "E15"; Value = "3.21%"}
    @{Cell = "D16"; Value = "3.025"}
    @{Cell = "E16"; Value = "0.43%"}
    @{Cell = "E17"; Value = "-1.74%"}
    @{Cell = "E18"; Value = "-2.39%"}
    @{Cell = "D19"; Value = "0.1331"}
    @{Cell = "E19"; Value = "-0.43%"}
    @{Cell = "D20"; Value = "0.03179"}
    @{Cell = "E20"; Value = "-1.84%"}
    @{Cell = "E21"; Value = "-1.94%"}
    @{Cell = "D22"; Value = "3.758"}
    @{Cell = "E22"; Value = "0.57%"}
    @{Cell = "D23"; Value = "0.04658"}
    @{Cell = "E23"; Value = "-0.76%"}
    @{Cell = "D24"; Value = "0.1373"}
    @{Cell = "E24"; Value = "2.48%"}
    @{Cell = "D25"; Value = "0.001243"}
    @{Cell = "E25"; Value = "0.19%"}
    @{Cell = "D26"; Value = "0.004549"}
    @{Cell = "E26"; Value = "5.91%"}
    @{Cell = "D27"; Value = "0.00009594"}
    @{Cell = "E27"; Value = "-1.03%"}
    @{Cell = "D28"; Value = "0.0001937"}
    @{Cell = "D40"; Value = "0.03650"}
    @{Cell = "E40"; Value = "-0.43%"}
    @{Cell = "D41"; Value = "0.1357"}
    @{Cell = "E41"; Value = "28.56%"}
    @{Cell = "D42"; Value = "0.006136"}
    @{Cell = "E42"; Value = "-1.12%"}
    @{Cell = "D43"; Value = "0.002638"}
    @{Cell = "E43"; Value = "1.54%"}
    @{Cell = "D44"; Value = "0.007978"}
    @{Cell = "E44"; Value = "-2.81%"}
    @{Cell = "D45"; Value = "0.00005381"}
    @{Cell = "E45"; Value = "1.70%"}
    @{Cell = "D47"; Value = "0.1449"}
    @{Cell = "E47"; Value = "-19.43%"}
    @{Cell = "D48"; Value = "0.002401"}
    @{Cell = "E48"; Value = "19.18%"}
    @{Cell = "D49"; Value = "0.00002099"}
    @{Cell = "D50"; Value = "0.0001999"}
)

foreach ($entry in $textCells) {
    $rng = $ws.Range($entry.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $entry.Value
}

# B/C columns are ordinary (non-numeric-looking) text, so a direct .Value
# assignment is safe and keeps the original (default) cell style untouched.
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("B42").Value = "KickToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
